$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 (copy formatting/style from H1, the existing header cell)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data columns I (I0) and J (IF) for rows 2-17
$iValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,5)
$jValues = @(6,7,6,4,5,7,6,6,6,4,6,6,5,5,5,6)

for ($r = 2; $r -le 17; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
